$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 2 (G=5489)
$ws.Cells.Item(2, 8).Value = 233
$ws.Cells.Item(2, 9).Value = 150
$ws.Cells.Item(2, 10).Value = 399
$ws.Cells.Item(2, 11).Value = 150
$ws.Cells.Item(2, 12).Value = 399
$ws.Cells.Item(2, 13).Value = -37
$ws.Cells.Item(2, 14).Value = -625
# Row 4 (G=5470)
$ws.Cells.Item(4, 8).Value = 94
$ws.Cells.Item(4, 9).Value = 90
$ws.Cells.Item(4, 10).Value = 102
$ws.Cells.Item(4, 11).Value = 90
$ws.Cells.Item(4, 12).Value = 102
$ws.Cells.Item(4, 13).Value = 24
$ws.Cells.Item(4, 14).Value = -330
# Row 9 (G=5487)
$ws.Cells.Item(9, 8).Value = 41.733334
$ws.Cells.Item(9, 9).Value = 35.090908
$ws.Cells.Item(9, 10).Value = 60
$ws.Cells.Item(9, 11).Value = 35.090908
$ws.Cells.Item(9, 12).Value = 60
$ws.Cells.Item(9, 13).Value = 133.909092
$ws.Cells.Item(9, 14).Value = -398
# Row 12 (G=5515)
$ws.Cells.Item(12, 8).Value = 274.66666
$ws.Cells.Item(12, 9).Value = 289.4
$ws.Cells.Item(12, 10).Value = 201
$ws.Cells.Item(12, 11).Value = 289.4
$ws.Cells.Item(12, 12).Value = 201
$ws.Cells.Item(12, 13).Value = -119.4
$ws.Cells.Item(12, 14).Value = -541
# Row 17 (G=38956)
$ws.Cells.Item(17, 8).Value = 1267119
$ws.Cells.Item(17, 10).Value = 1389598.2
$ws.Cells.Item(17, 12).Value = 4168794.6
$ws.Cells.Item(17, 14).Value = -4169130.6
# Row 41 (G=5478)
$ws.Cells.Item(41, 8).Value = 11111582
$ws.Cells.Item(41, 9).Value = 12346180
$ws.Cells.Item(41, 10).Value = 200
$ws.Cells.Item(41, 11).Value = 12346180
$ws.Cells.Item(41, 12).Value = 200
$ws.Cells.Item(41, 13).Value = -12345740
$ws.Cells.Item(41, 14).Value = -1080
# Row 47 (G=2169)
$ws.Cells.Item(47, 8).Value = 1566.75
$ws.Cells.Item(47, 9).Value = 1566.75
$ws.Cells.Item(47, 11).Value = 1566.75
$ws.Cells.Item(47, 13).Value = -594.75
# Row 95 (G=18200)
$ws.Cells.Item(95, 8).Value = 500624
$ws.Cells.Item(95, 10).Value = 500624
$ws.Cells.Item(95, 12).Value = 500624
$ws.Cells.Item(95, 14).Value = -506116
# Row 96 (G=19894)
$ws.Cells.Item(96, 8).Value = 750.4
$ws.Cells.Item(96, 9).Value = 563
$ws.Cells.Item(96, 10).Value = 1500
$ws.Cells.Item(96, 11).Value = 1689
$ws.Cells.Item(96, 12).Value = 4500
$ws.Cells.Item(96, 13).Value = -316
$ws.Cells.Item(96, 14).Value = -7246
# Row 113 (G=27775)
$ws.Cells.Item(113, 8).Value = 1891.3334
$ws.Cells.Item(113, 9).Value = 1891.3334
$ws.Cells.Item(113, 11).Value = 1891.3334
$ws.Cells.Item(113, 13).Value = 1362.6666
# Row 137 (G=44013)
$ws.Cells.Item(137, 8).Value = 90910610
$ws.Cells.Item(137, 9).Value = 125001160
$ws.Cells.Item(137, 10).Value = 2468.6667
$ws.Cells.Item(137, 11).Value = 375003480
$ws.Cells.Item(137, 12).Value = 7406.000100000001
$ws.Cells.Item(137, 13).Value = -375000930
$ws.Cells.Item(137, 14).Value = -12506.0001

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 32 (G=44147)
$ws.Cells.Item(32, 8).Value = 15486.454
$ws.Cells.Item(32, 9).Value = 2266.3562
$ws.Cells.Item(32, 11).Value = 2266.3562
$ws.Cells.Item(32, 13).Value = -1979.3562
# Row 92 (G=18050)
$ws.Cells.Item(92, 8).Value = 33399.5
$ws.Cells.Item(92, 10).Value = 33399.5
$ws.Cells.Item(92, 12).Value = 33399.5
$ws.Cells.Item(92, 14).Value = -38391.5
# Row 95 (G=18204)
$ws.Cells.Item(95, 8).Value = 116069.336
$ws.Cells.Item(95, 10).Value = 116069.336
$ws.Cells.Item(95, 12).Value = 116069.336
$ws.Cells.Item(95, 14).Value = -121561.336

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 4 (G=3742)
$ws.Cells.Item(4, 8).Value = 500000500
$ws.Cells.Item(4, 10).Value = 500000500
$ws.Cells.Item(4, 12).Value = 500000500
$ws.Cells.Item(4, 14).Value = -500000724
# Row 22 (G=5367)
$ws.Cells.Item(22, 8).Value = 470.93332
$ws.Cells.Item(22, 9).Value = 296.4
$ws.Cells.Item(22, 10).Value = 820
$ws.Cells.Item(22, 11).Value = 296.4
$ws.Cells.Item(22, 12).Value = 820
$ws.Cells.Item(22, 13).Value = 53.60000000000002
$ws.Cells.Item(22, 14).Value = -1520
# Row 31 (G=44023)
$ws.Cells.Item(31, 8).Value = 2022.8125
$ws.Cells.Item(31, 9).Value = 1195.1305
$ws.Cells.Item(31, 10).Value = 4138
$ws.Cells.Item(31, 11).Value = 1195.1305
$ws.Cells.Item(31, 12).Value = 4138
$ws.Cells.Item(31, 13).Value = -900.1305
$ws.Cells.Item(31, 14).Value = -4728
# Row 34 (G=44023)
$ws.Cells.Item(34, 8).Value = 2022.8125
$ws.Cells.Item(34, 9).Value = 1195.1305
$ws.Cells.Item(34, 10).Value = 4138
$ws.Cells.Item(34, 11).Value = 1195.1305
$ws.Cells.Item(34, 12).Value = 4138
$ws.Cells.Item(34, 13).Value = -993.1305
$ws.Cells.Item(34, 14).Value = -4542

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 33 (G=4867)
$ws.Cells.Item(33, 8).Value = 49.066666
$ws.Cells.Item(33, 9).Value = 48.22222
$ws.Cells.Item(33, 10).Value = 50.333332
$ws.Cells.Item(33, 11).Value = 289.33332
$ws.Cells.Item(33, 12).Value = 301.999992
$ws.Cells.Item(33, 13).Value = -6.333320000000015
$ws.Cells.Item(33, 14).Value = -867.999992
# Row 75 (G=12863)
$ws.Cells.Item(75, 8).Value = 1000
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 12).Value = 0
$ws.Cells.Item(75, 14).ClearContents()
# Row 78 (G=12863)
$ws.Cells.Item(78, 8).Value = 1000
$ws.Cells.Item(78, 10).Value = 0
$ws.Cells.Item(78, 12).Value = 0
$ws.Cells.Item(78, 14).ClearContents()
# Row 131 (G=36060)
$ws.Cells.Item(131, 8).Value = 1335.9688
$ws.Cells.Item(131, 9).Value = 625
$ws.Cells.Item(131, 10).Value = 1358.9032
$ws.Cells.Item(131, 11).Value = 1875
$ws.Cells.Item(131, 12).Value = 4076.7096
$ws.Cells.Item(131, 13).Value = 3165
$ws.Cells.Item(131, 14).Value = -14156.7096
# Row 132 (G=43972)
$ws.Cells.Item(132, 8).Value = 1103.2222
$ws.Cells.Item(132, 9).Value = 0
$ws.Cells.Item(132, 10).Value = 1103.2222
$ws.Cells.Item(132, 11).Value = 0
$ws.Cells.Item(132, 12).ClearContents()
$ws.Cells.Item(132, 13).Value = 9928.9998
$ws.Cells.Item(132, 14).Value = -14988.9998
# Row 137 (G=44088)
$ws.Cells.Item(137, 8).Value = 3611000
$ws.Cells.Item(137, 9).Value = 6670477
$ws.Cells.Item(137, 10).Value = 80834.46000000001
$ws.Cells.Item(137, 11).Value = 20011431
$ws.Cells.Item(137, 12).Value = 242503.38
$ws.Cells.Item(137, 13).Value = -20006331
$ws.Cells.Item(137, 14).Value = -252703.38

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 5 (G=1681)
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 10).Value = 0
$ws.Cells.Item(5, 12).ClearContents()
$ws.Cells.Item(5, 14).Value = 0
# Row 52 (G=4147)
$ws.Cells.Item(52, 8).Value = 44000
$ws.Cells.Item(52, 10).Value = 44000
$ws.Cells.Item(52, 12).Value = 44000
$ws.Cells.Item(52, 14).Value = -44518
# Row 133 (G=41854)
$ws.Cells.Item(133, 8).Value = 19120.715
$ws.Cells.Item(133, 10).Value = 19120.715
$ws.Cells.Item(133, 12).Value = 19120.715
$ws.Cells.Item(133, 14).Value = -29240.715

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 2 (G=2631)
$ws.Cells.Item(2, 8).Value = 0
$ws.Cells.Item(2, 10).Value = 0
$ws.Cells.Item(2, 12).ClearContents()
$ws.Cells.Item(2, 14).Value = 0
# Row 7 (G=36249)
$ws.Cells.Item(7, 8).Value = 3523.5293
$ws.Cells.Item(7, 9).Value = 4000
$ws.Cells.Item(7, 10).Value = 3460
$ws.Cells.Item(7, 11).Value = 4000
$ws.Cells.Item(7, 12).Value = 3460
$ws.Cells.Item(7, 13).Value = -3888
$ws.Cells.Item(7, 14).Value = -3684
# Row 46 (G=5282)
$ws.Cells.Item(46, 8).Value = 750.6429000000001
$ws.Cells.Item(46, 9).Value = 649.8333
$ws.Cells.Item(46, 10).Value = 826.25
$ws.Cells.Item(46, 11).Value = 649.8333
$ws.Cells.Item(46, 12).Value = 826.25
$ws.Cells.Item(46, 13).Value = -461.8333
$ws.Cells.Item(46, 14).Value = -1202.25
# Row 126 (G=36249)
$ws.Cells.Item(126, 8).Value = 3523.5293
$ws.Cells.Item(126, 9).Value = 4000
$ws.Cells.Item(126, 10).Value = 3460
$ws.Cells.Item(126, 11).Value = 12000
$ws.Cells.Item(126, 12).Value = 10380
$ws.Cells.Item(126, 13).Value = -9530
$ws.Cells.Item(126, 14).Value = -15320
# Row 132 (G=44058)
$ws.Cells.Item(132, 8).Value = 3871.762
$ws.Cells.Item(132, 9).Value = 2282.077
$ws.Cells.Item(132, 11).Value = 6846.231000000001
$ws.Cells.Item(132, 13).Value = -4316.231000000001
# Row 136 (G=44060)
$ws.Cells.Item(136, 8).Value = 8137.5264
$ws.Cells.Item(136, 9).Value = 5952.909
$ws.Cells.Item(136, 10).Value = 11141.375
$ws.Cells.Item(136, 11).Value = 17858.727
$ws.Cells.Item(136, 12).Value = 33424.125
$ws.Cells.Item(136, 13).Value = -15308.727
$ws.Cells.Item(136, 14).Value = -38524.125

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 74 (G=19022)
$ws.Cells.Item(74, 8).Value = 13276.5
$ws.Cells.Item(74, 10).Value = 12940.5
$ws.Cells.Item(74, 12).Value = 12940.5
$ws.Cells.Item(74, 14).Value = -14812.5
# Row 77 (G=19022)
$ws.Cells.Item(77, 8).Value = 13276.5
$ws.Cells.Item(77, 10).Value = 12940.5
$ws.Cells.Item(77, 12).Value = 38821.5
$ws.Cells.Item(77, 14).Value = -48181.5
# Row 132 (G=44029)
$ws.Cells.Item(132, 8).Value = 3497.52
$ws.Cells.Item(132, 9).Value = 3544.375
$ws.Cells.Item(132, 10).Value = 3414.2222
$ws.Cells.Item(132, 11).Value = 10633.125
$ws.Cells.Item(132, 12).Value = 10242.6666
$ws.Cells.Item(132, 13).Value = -8103.125
$ws.Cells.Item(132, 14).Value = -15302.6666
